$wb = $excel.ActiveWorkbook

# --- Sheet "Players": insert a new row at row 18 (shifts old rows 18-42 down to 19-43) ---
$ws = $wb.Worksheets.Item("Players")
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new game log entry.
# Column A holds dates as plain text (matching the rest of the sheet), so force
# a text number format on that single cell before writing the date-shaped string
# to keep Excel from auto-converting it into a date serial number.
$ws.Cells.Item(18,1).NumberFormat = "@"
$ws.Cells.Item(18,1).Value  = "2026-02-04"
$ws.Cells.Item(18,2).Value  = "The Oddities"
$ws.Cells.Item(18,3).Value  = "Yes"
$ws.Cells.Item(18,4).Value  = "Rashaun Agee"
$ws.Cells.Item(18,5).Value  = "TA&M"
$ws.Cells.Item(18,6).Value  = "TA&M@ALA"
$ws.Cells.Item(18,7).Value  = "Wed, February 4th at 7:00 PM EST"
$ws.Cells.Item(18,8).Value  = -1
$ws.Cells.Item(18,9).Value  = 0
$ws.Cells.Item(18,10).Value = 0
$ws.Cells.Item(18,11).Value = 0
$ws.Cells.Item(18,12).Value = 0
$ws.Cells.Item(18,13).Value = 0
$ws.Cells.Item(18,14).Value = 0
$ws.Cells.Item(18,15).Value = 0
$ws.Cells.Item(18,16).Value = 0
$ws.Cells.Item(18,17).Value = 0
$ws.Cells.Item(18,18).Value = 1
$ws.Cells.Item(18,19).Value = 0
$ws.Cells.Item(18,20).Value = 0
$ws.Cells.Item(18,21).Value = 0
$ws.Cells.Item(18,22).Value = 0

# Widen column G ("status") from 8 to 34 characters.
# ColumnWidth is offset from the stored OOXML width by ~0.8333 (5/6) chars,
# so subtract that to land exactly on 34.
$ws.Columns.Item(7).ColumnWidth = 33.166666666666664

# --- Sheet "OwnerTotals": update "The Oddities" totals (row 7) ---
$ws2 = $wb.Worksheets.Item("OwnerTotals")
$ws2.Cells.Item(7,2).Value = 9
$ws2.Cells.Item(7,3).Value = 2
